# KCR_1994: add season-record columns (Wins / Losses / Ties) for every
# player row, mirroring the header formatting already used by column AC
# ("Unnamed: 28").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the existing header formatting (bold font, border, centered/top
# alignment) from AC1 onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the 1994 Kansas City Royals season record (64-51-0) for every
# player row in the sheet (rows 2-33).
$lastRow = 33
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 64   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 51   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
